$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Alfred Pennyworth's photo set was re-uploaded (ap4 swapped to .png) - update first
$ws.Range("E11").Value = "ap1.jpg,ap2.jpg,ap3.jpg,ap4.png,ap5.jpg"

# Fill in / refresh the rest of the photo sets in row order
$ws.Range("E2").Value  = "ss1.jpg,ss2.jpg,ss3.jpg,ss4.jpg,ss5.jpg"
$ws.Range("E3").Value  = "ph1.jpg,ph2,jpg,ph3,jpg,ph4.jpg,ph5.jpg"
$ws.Range("E4").Value  = "zw1.jpg,zw2.jpg,zw3.jpg,zw4.jpg,zw5.jpg"
$ws.Range("E5").Value  = "mf1.jpg,mf2.jpg,mf3.jpg,mf4.jpg,mf5.jpg"
$ws.Range("E7").Value  = "rt1.jpg,rt2.jpg,rt3.jpg,rt4.jpg,rt5.jpg"
$ws.Range("E9").Value  = "mff1.jpg,mff2.jpg,mff3.jpg,mff4.jpg,mff5.jpg"
$ws.Range("E6").Value  = "bb1.jpg,bb2.jpg,bb3.jpg,bb4.jpg,bb5.jpg"
$ws.Range("E8").Value  = "pf1.jpg,pf2.jpg,pf3.jpg,pf4.jpg,pf5.jpg"
$ws.Range("E10").Value = "jm1.jpg,jm2.jpg,jm3.jpg,jm4.jpg,jm5.jpg"

$ws.Range("E13").Select()

$wb.Save()
